# Adds two new expense rows (2 and 3) to the "Expense" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: # , Date and Time, Amount, Currency, Lender Name, Borrower Name
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "2022-01-13 16:50:17.078917"
$ws.Range("C2").Value = 500
$ws.Range("D2").Value = "EUR"
$ws.Range("E2").Value = "qcvcvbcvy"
$ws.Range("F2").Value = "qcvcvbcvy"

# Row 3
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "2022-01-13 16:50:26.841231"
$ws.Range("C3").Value = 500
$ws.Range("D3").Value = "EUR"
$ws.Range("E3").Value = "qcvcvbcvy"
$ws.Range("F3").Value = "qcvcvbcvy"
